# Implement database schema migration (issue #31): append the new
# row (row 37) of captured frame data to each of the four sheets,
# mirroring the previous row (row 36) but stamped with the next day's
# timestamp.

$wb = $excel.ActiveWorkbook

$newRow = 37
$prevRow = 36

# Per-sheet payload for the new row: B/C/D/E are text-ish fields (hex byte
# dumps), F/G/H/I are numeric decodes. A is the new timestamp (shared by
# all four sheets).
$sheetsData = @{
    "DE_LFT_#1" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x70"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 368
        I = 14
    }
    "DE_LFT_#2" = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x70"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 368
        I = 14
    }
    "DE_PLT_#1" = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x80"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 128
        I = 7
    }
    "DE_PLT_#2" = @{
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x80"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 128
        I = 3
    }
}

foreach ($sheetName in $sheetsData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetsData[$sheetName]

    # Column A: new timestamp, same numeric format/style as the row above.
    $aCell = $ws.Cells.Item($newRow, 1)
    $aCell.Value = 45823.43395833333
    $aCell.NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
